$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) column stores values as text in the source data (e.g.
# thousand-separated prices like "24.270.10"). For the subset of new
# values that look like plain decimal numbers (e.g. "1.000", "49.50"),
# Excel would otherwise auto-convert them to numeric values and drop
# significant trailing/leading zeros, so mark those cells as Text first
# so the exact original string is preserved.
$priceCellsNeedingTextFormat = @("D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D16", "D18", "D19", "D20", "D21", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $priceCellsNeedingTextFormat) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "24.270.10"
$ws.Range("E2").Value = "  -2.82%  "
$ws.Range("D3").Value = "1.649.90"
$ws.Range("E3").Value = "  -3.39%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "309.61"
$ws.Range("E5").Value = "  -2.09%  "
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "0.3895"
$ws.Range("E7").Value = "  -1.48%  "
$ws.Range("D8").Value = "0.3889"
$ws.Range("E8").Value = "  -3.66%  "
$ws.Range("D9").Value = "1.002"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("D10").Value = "1.374"
$ws.Range("E10").Value = "  -7.53%  "
$ws.Range("D11").Value = "49.50"
$ws.Range("E11").Value = "  -6.20%  "
$ws.Range("D12").Value = "0.08501"
$ws.Range("E12").Value = "  -3.60%  "
$ws.Range("E13").Value = "  -5.72%  "
$ws.Range("D14").Value = "7.178"
$ws.Range("E14").Value = "  -4.17%  "
$ws.Range("E15").Value = "  -4.79%  "
$ws.Range("D16").Value = "7.534"
$ws.Range("E16").Value = "  -5.81%  "
$ws.Range("D17").Value = "1.649.42"
$ws.Range("E17").Value = "  -5.04%  "
$ws.Range("D18").Value = "94.95"
$ws.Range("E18").Value = "  -1.38%  "
$ws.Range("D19").Value = "21.16"
$ws.Range("E19").Value = "  +2.58%  "
$ws.Range("D20").Value = "0.06899"
$ws.Range("E20").Value = "  -3.94%  "
$ws.Range("D21").Value = "6.972"
$ws.Range("E21").Value = "  -5.34%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "13.88"
$ws.Range("D24").Value = "24.272.90"
$ws.Range("E24").Value = "  -2.78%  "
$ws.Range("D25").Value = "2.370"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").Value = "2.771"
$ws.Range("E26").Value = "  -7.02%  "
$ws.Range("D27").Value = "22.63"
$ws.Range("E27").Value = "  -4.60%  "
$ws.Range("D28").Value = "158.43"
$ws.Range("E28").Value = "  -1.99%  "
$ws.Range("D29").Value = "8.591"
$ws.Range("E29").Value = "  +1.44%  "
$ws.Range("D30").Value = "143.30"
$ws.Range("E30").Value = "  -4.78%  "
$ws.Range("D31").Value = "5.368"
$ws.Range("E31").Value = "  -13.70%  "
$ws.Range("D32").Value = "2.459"
$ws.Range("E32").Value = "  -4.70%  "
$ws.Range("D33").Value = "1.829.50"
$ws.Range("E33").Value = "  -3.82%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "6.991"
$ws.Range("E34").Value = "  -2.67%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.08168"
$ws.Range("E35").Value = "  -4.49%  "
$ws.Range("D36").Value = "0.9965"
$ws.Range("E36").Value = "  -4.89%  "
$ws.Range("D37").Value = "0.02955"
$ws.Range("E37").Value = "  -5.71%  "
$ws.Range("D38").Value = "0.2730"
$ws.Range("E38").Value = "  -4.51%  "
$ws.Range("D39").Value = "0.09338"
$ws.Range("E39").Value = "  -1.96%  "
$ws.Range("D40").Value = "1.485"
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").Value = "10.06"
$ws.Range("E41").Value = "  -7.21%  "
$ws.Range("D42").Value = "0.7682"
$ws.Range("E42").Value = "  -7.17%  "
$ws.Range("D43").Value = "13.20"
$ws.Range("E43").Value = "  -5.88%  "
$ws.Range("D44").Value = "16.05"
$ws.Range("E44").Value = "  -8.28%  "
$ws.Range("E45").Value = "  -6.48%  "
$ws.Range("D46").Value = "0.6923"
$ws.Range("E46").Value = "  -6.47%  "
$ws.Range("D47").Value = "4.108"
$ws.Range("E47").Value = "  -3.24%  "
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").Value = "0.08467"
$ws.Range("E49").Value = "  -3.07%  "
$ws.Range("D50").Value = "1.276"
$ws.Range("E50").Value = "  -7.42%  "
$ws.Range("D51").Value = "134.87"
$ws.Range("E51").Value = "  -3.09%  "
